$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates (odds refreshed) ---
$ws.Range("N2").Value = 3.15
$ws.Range("O2").Value = 1.29
$ws.Range("Q2").Value = 1.86
$ws.Range("X2").Value = 16.5
$ws.Range("AC2").Value = 10.5
$ws.Range("AD2").Value = 26
$ws.Range("AF2").Value = 11.5
$ws.Range("AH2").Value = 26
$ws.Range("AJ2").Value = 20
$ws.Range("AK2").Value = 21
$ws.Range("AN2").Value = 12

# --- Row 3 updates (York City vs Rochdale odds refreshed) ---
$ws.Range("F3").Value = 1.85
$ws.Range("G3").Value = 2.44
$ws.Range("H3").Value = 3.15
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 3.7
$ws.Range("K3").Value = 950
$ws.Range("L3").Value = 1.24
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 2.02
$ws.Range("O3").Value = 1.22
$ws.Range("P3").Value = 2.02
$ws.Range("Q3").Value = 1.6
$ws.Range("R3").Value = 1.28
$ws.Range("S3").Value = 2
$ws.Range("T3").Value = 1.03
$ws.Range("U3").Value = 1.03
$ws.Range("V3").Value = 1.28
$ws.Range("W3").Value = 1.7
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# --- Row 4 (new game added: Colombian Primera B) ---
$ws.Range("A4").Value = "Colombian Primera B"

# Force B4/C4 to stay literal text (not auto-converted to date/time serials),
# then drop back to the default "Normal" style so no stray cell formatting
# is left behind, matching League/Home/Away neighbour cells which have no
# explicit style either.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2025-11-11"
$ws.Range("B4").Style = "Normal"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "18:00:00"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "Cucuta Deportivo"
$ws.Range("E4").Value = "Jaguares de Cordoba"
$ws.Range("F4").Value = 2.14
$ws.Range("G4").Value = 2.98
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 5.1
$ws.Range("J4").Value = 2.46
$ws.Range("K4").Value = 4.1
$ws.Range("L4").Value = 1.01
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 1.9
$ws.Range("O4").Value = 1.01
$ws.Range("P4").Value = 1.44
$ws.Range("Q4").Value = 2.2
$ws.Range("R4").Value = 1.14
$ws.Range("S4").Value = 2.2
$ws.Range("T4").Value = 1.04
$ws.Range("U4").Value = 1.04
$ws.Range("V4").Value = 1.24
$ws.Range("W4").Value = 1.5
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 1000
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 1000
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000
